$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update email, clear date (backend now only appends email without date on this row)
$ws.Range("A2").Value = "as@gmail.com"
$ws.Range("B2").Value = $null

# Row 3: update email and date
$ws.Range("A3").Value = "yadav@gmail.com"
$ws.Range("B3").Value = "11/1/2025, 8:50:26 PM"
